$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to remain Text storage even when the new value
# looks like a plain number (e.g. "1.00", "39.70"), while preserving the
# cell's original style index (round-trip Style so no style/format diff).
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "42.980.68"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.299.55"
$ws.Range("E3").Value = "  -0.18%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue $ws.Range("D5") "312.49"
$ws.Range("E5").Value = "  -3.10%  "
Set-TextValue $ws.Range("D6") "104.64"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -1.45%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -0.77%  "
Set-TextValue $ws.Range("D10") "39.70"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  -1.16%  "
Set-TextValue $ws.Range("D12") "8.26"
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("E14").Value = "  +0.81%  "
Set-TextValue $ws.Range("D15") "15.34"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").Value = "2.646.10"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "2.297.09"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "42.819.67"
$ws.Range("E18").Value = "  -0.14%  "
Set-TextValue $ws.Range("D19") "7.32"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("E20").Value = "  -1.82%  "
Set-TextValue $ws.Range("D21") "13.44"
$ws.Range("E21").Value = "  +1.28%  "
Set-TextValue $ws.Range("D22") "73.41"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D23") "3.43"
$ws.Range("E23").Value = "  -5.34%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D24") "267.45"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("E26").Value = "  +0.42%  "
Set-TextValue $ws.Range("D27") "10.82"
$ws.Range("E27").Value = "  -1.46%  "
Set-TextValue $ws.Range("D28") "7.11"
$ws.Range("E28").Value = "  +15.04%  "
Set-TextValue $ws.Range("D29") "2.29"
$ws.Range("E29").Value = "  -1.42%  "
Set-TextValue $ws.Range("D30") "22.33"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  -5.29%  "
Set-TextValue $ws.Range("D32") "164.97"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D34") "2.63"
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D35") "0.130"
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("E36").Value = "  -2.90%  "
Set-TextValue $ws.Range("D37") "4.56"
$ws.Range("E37").Value = "  -2.17%  "
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  -2.73%  "
Set-TextValue $ws.Range("D41") "108.06"
$ws.Range("E41").Value = "  +5.24%  "
$ws.Range("E42").Value = "  +0.68%  "
Set-TextValue $ws.Range("D43") "71.11"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  +0.77%  "
Set-TextValue $ws.Range("D45") "1.01"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D46") "12.15"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.742.96"
$ws.Range("E47").Value = "  +9.58%  "
Set-TextValue $ws.Range("D48") "110.56"
$ws.Range("E48").Value = "  -3.63%  "
Set-TextValue $ws.Range("D49") "77.62"
$ws.Range("E49").Value = "  -6.41%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D50") "5.16"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D51") "8.65"
$ws.Range("E51").Value = "  -3.04%  "
